$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2, 3 and 4 are observation records sharing the same location/date/
# observer, but differing in species-specific data (columns A, B, D, E, F,
# G, H, Q, R, Z, AB and the optional public-comment column AC). This edit
# rotates those values one row down, wrapping row 4 back into row 2:
#   new row2 <- old row3
#   new row3 <- old row4
#   new row4 <- old row2

$cols = @("A","B","D","E","F","G","H","Q","R","Z","AB","AC")

# Capture the current ("before") values for rows 2, 3 and 4.
$old2 = @{}
$old3 = @{}
$old4 = @{}
foreach ($c in $cols) {
    $old2[$c] = $ws.Range("${c}2").Value2
    $old3[$c] = $ws.Range("${c}3").Value2
    $old4[$c] = $ws.Range("${c}4").Value2
}

# Write rotated values back: row2 = old row3, row3 = old row4, row4 = old row2
foreach ($c in $cols) {
    $ws.Range("${c}2").Value2 = $old3[$c]
    $ws.Range("${c}3").Value2 = $old4[$c]
    $ws.Range("${c}4").Value2 = $old2[$c]
}
